# This script swaps the data between rows 2 & 3 (full record swap for
# columns A, I, J, K, Q, R, AC) and rows 4 & 5 (swap of A, Q, R only),
# matching the upstream dataset re-ordering/correction described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap full rows 2 and 3 (columns A, I, J, K, Q, R, AC) ---
# NOTE: use .Value2 for reads — .Value's getter in this runtime returns a
# placeholder descriptor string instead of the real cell content.

$row2_A = $ws.Range("A2").Value2
$row2_I = $ws.Range("I2").Value2
$row2_J = $ws.Range("J2").Value2
$row2_K = $ws.Range("K2").Value2
$row2_Q = $ws.Range("Q2").Value2
$row2_R = $ws.Range("R2").Value2
$row2_AC = $ws.Range("AC2").Value2

$row3_A = $ws.Range("A3").Value2
$row3_I = $ws.Range("I3").Value2
$row3_J = $ws.Range("J3").Value2
$row3_K = $ws.Range("K3").Value2
$row3_Q = $ws.Range("Q3").Value2
$row3_R = $ws.Range("R3").Value2
$row3_AC = $ws.Range("AC3").Value2

$ws.Range("A2").Value = $row3_A
$ws.Range("I2").Value = $row3_I
$ws.Range("J2").Value = $row3_J
$ws.Range("K2").Value = $row3_K
$ws.Range("Q2").Value = $row3_Q
$ws.Range("R2").Value = $row3_R
$ws.Range("AC2").Value = $row3_AC

$ws.Range("A3").Value = $row2_A
$ws.Range("I3").Value = $row2_I
$ws.Range("J3").Value = $row2_J
$ws.Range("K3").Value = $row2_K
$ws.Range("Q3").Value = $row2_Q
$ws.Range("R3").Value = $row2_R
$ws.Range("AC3").Value = $row2_AC

# --- Swap rows 4 and 5 (columns A, Q, R only) ---

$row4_A = $ws.Range("A4").Value2
$row4_Q = $ws.Range("Q4").Value2
$row4_R = $ws.Range("R4").Value2

$row5_A = $ws.Range("A5").Value2
$row5_Q = $ws.Range("Q5").Value2
$row5_R = $ws.Range("R5").Value2

$ws.Range("A4").Value = $row5_A
$ws.Range("Q4").Value = $row5_Q
$ws.Range("R4").Value = $row5_R

$ws.Range("A5").Value = $row4_A
$ws.Range("Q5").Value = $row4_Q
$ws.Range("R5").Value = $row4_R
